$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fall 22 week 11 complete plus 9 ball skill level evals
# Update matchup results: C3 (2/1 -> 0/3) and B5 (0/2 -> 3/0)
$ws.Range("C3").Value = "0/3"
$ws.Range("B5").Value = "3/0"
